# feat: add 2022-Q1 data
#
# Before: sheets = [2021-Q2, 2021-Q3, 总计]
# After:  sheets = [2021-Q2, 2021-Q3, 2022-Q1, 总计]
#   - "总计" (index 3) is renamed to "2022-Q1" and its content is replaced
#     with the new quarter's per-fund holdings table.
#   - a fresh "总计" sheet is inserted after it (a copy of the old "总计"
#     sheet, so it keeps the same look/formatting) with a new summary row
#     for 2022-Q1 added on top of the existing history.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(3)

# Duplicate the existing "总计" sheet right after itself - the copy keeps
# all formatting/styles and will become the new "总计" sheet, while the
# original gets repurposed below into "2022-Q1".
$totalSheet.Copy($null, $totalSheet)
$newTotalSheet = $wb.Worksheets.Item(4)

$totalSheet.Name = "2022-Q1"
$newTotalSheet.Name = "总计"

# ---------------------------------------------------------------------
# 1) "2022-Q1" sheet: replace the old 总计-style table with the new
#    per-fund holdings table (basé on the other quarter sheets' layout).
# ---------------------------------------------------------------------
$ws = $totalSheet

# Extend the header styling (s=2, same as B1:D1) across the new columns.
$ws.Range("D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)

# Extend column-A's style (s=2) down to the new 3rd data row.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# B (基金代码) and D:G are textual in the source data (fund codes keep
# leading context / the numbers keep their original string formatting,
# e.g. trailing zeros), so force text before assigning, then drop back
# to the "Normal" style afterwards so no stray number-format style is
# left attached to the cell (matches the other quarter sheets, where
# these data cells carry no explicit style).
$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("D2:G4").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "513050"
$ws.Range("C2").Value = "易方达中证海外中国互联网50 QDII-ETF"
$ws.Range("D2").Value = "350.10"
$ws.Range("E2").Value = "98.05"
$ws.Range("F2").Value = "3.66"
$ws.Range("G2").Value = "12.8137"
$ws.Range("H2").Value = 6

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "159605"
$ws.Range("C3").Value = "广发中证海外中国互联网30（QDII-ETF）"
$ws.Range("D3").Value = "29.04"
$ws.Range("E3").Value = "98.61"
$ws.Range("F3").Value = "6.82"
$ws.Range("G3").Value = "1.9805"
$ws.Range("H3").Value = 6

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "159607"
$ws.Range("C4").Value = "嘉实中证海外中国互联网30ETF（QDII）"
$ws.Range("D4").Value = "5.79"
$ws.Range("E4").Value = "98.25"
$ws.Range("F4").Value = "6.87"
$ws.Range("G4").Value = "0.3978"
$ws.Range("H4").Value = 6

$ws.Range("B2:B4").Style = "Normal"
$ws.Range("D2:G4").Style = "Normal"

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert the 2022-Q1 summary row on top, pushing the
#    existing 2021-Q3 / 2021-Q2 rows down by one.
# ---------------------------------------------------------------------
$ts = $newTotalSheet

# Make room for the new row, keeping the old row 2/3 formatting in place.
$ts.Range("A3:D3").Copy()
$ts.Range("A4:D4").PasteSpecial(-4122)
$ts.Range("A2:D2").Copy()
$ts.Range("A3:D3").PasteSpecial(-4122)

$ts.Range("A2").Value = 0
$ts.Range("B2").Value = "2022-Q1"
$ts.Range("C2").Value = 3
$ts.Range("D2").Value = 15.19

$ts.Range("A3").Value = 1
$ts.Range("B3").Value = "2021-Q3"
$ts.Range("C3").Value = 2
$ts.Range("D3").Value = 10.39

$ts.Range("A4").Value = 2
$ts.Range("B4").Value = "2021-Q2"
$ts.Range("C4").Value = 3
$ts.Range("D4").Value = 6.66

Write-Output "ok"
